$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $last)
$ws.Name = "Evaluation Warning"

$msg = "Evaluation Only. Created with Aspose.Cells for .NET.Copyright 2003 - 2020 Aspose Pty Ltd."
$cell = $ws.Range("A5")
$cell.Value = $msg
$ws.Rows.Item(5).RowHeight = 23.25

$cell.Font.Name = "Arial"
$cell.Font.Size = 18
$cell.Font.Bold = $true
$cell.Font.Italic = $true
$cell.Font.Color = 16711680
